$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename existing shared string used in C2 ("BreedingGrounds" -> "SandBox")
$ws.Range("C2").Value = "SandBox"

# Fill in the newly added rows 19-21 (B/C/D) that previously only had column A filled
$ws.Range("B19").Value = "Interruptor"
$ws.Range("C19").Value = "Entrance-BreedingGrounds"
$ws.Range("D19").Value = "Opens First Door before Light"

$ws.Range("B20").Value = "Interruptor"
$ws.Range("C20").Value = "Entrance-BreedingGrounds"
$ws.Range("D20").Value = "Opens Second Door before Light"

$ws.Range("B21").Value = "Interruptor"
$ws.Range("C21").Value = "Entrance-BreedingGrounds"
$ws.Range("D21").Value = "Opens Third Door before Light"

# Move view/selection similar to target state
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D22").Select()
